$wb = $excel.ActiveWorkbook

# --- 1) Rename header cells on the existing sheets -------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after "Monthly Trend" --------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# --- 3) Write the header row -------------------------------------------------
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Reuse the existing bold/bordered/centered header style (same as the other
# two sheets) rather than synthesising a brand-new one.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# --- 4) Write the forecast data rows ----------------------------------------
$data = @(
    @(45508.99999999999, 8, 0.7954742517349309, 15.95388935517973),
    @(45515.99999999999, 8, -0.4518208260067587, 15.76955866426955),
    @(45529.99999999999, 7, -0.6504731851627906, 14.05089570949595),
    @(45564.99999999999, 5, -2.630469695815297, 12.14530178961112),
    @(45613.99999999999, 2, -6.472520767281297, 8.425753738007957),
    @(45620.99999999999, 1, -6.377707741838797, 9.029930878234024),
    @(45627.99999999999, 1, -6.990414935327292, 8.333017529419404),
    @(45634.99999999999, 0, -7.464346367183849, 7.086450828425741),
    @(45641.99999999999, 0, -7.590510174291037, 7.221806830418761),
    @(45648.99999999999, 0, -8.201233606931616, 6.300375526523325),
    @(45655.99999999999, 0, -8.944274865370167, 5.998359603408207),
    @(45662.99999999999, 0, -9.510109854880227, 5.688761661503247),
    @(45669.99999999999, 0, -9.472000660268881, 5.227279288543895),
    @(45676.99999999999, 0, -9.710426799744113, 5.007441618911245)
)

$row = 2
foreach ($d in $data) {
    $ws3.Cells.Item($row, 1).Value = $d[0]
    $ws3.Cells.Item($row, 2).Value = $d[1]
    $ws3.Cells.Item($row, 3).Value = $d[2]
    $ws3.Cells.Item($row, 4).Value = $d[3]
    $row++
}

# Apply the same date/time number-format style used by column A on the other
# two sheets to the new "ds" column's data rows.
$ws1.Range("A2").Copy()
$ws3.Range("A2:A15").PasteSpecial(-4122)

$ws3.Range("A1").Select()
